# newsbot state update: append latest scraped headlines to "Historico"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 79
$ws.Range("A79").Value = '05/01/2026 12:29:53'
$ws.Range("B79").Value = '05/01 12:14'
$ws.Range("C79").Value = 'Metrópoles'
$ws.Range("D79").Value = 'STF invalida lei que mudava critério em concurso para juiz estadual'
$ws.Range("E79").Value = 'https://www.metropoles.com/brasil/stf-invalida-lei-que-mudava-criterio-em-concurso-para-juiz-estadual'
$ws.Range("F79").Value = 'stf'
$ws.Range("G79").Value = 'Lei de MT estabelecia idade mínima para inscrição em concurso da magistratura. O STF entendeu que estados e DF não têm essa competência'

# Row 80 (no Score/keyword, no snippet)
$ws.Range("A80").Value = '05/01/2026 12:29:53'
$ws.Range("B80").Value = '05/01 12:13'
$ws.Range("C80").Value = 'Metrópoles'
$ws.Range("D80").Value = 'Em despacho, ministro do TCU prevê cautelar contra BC no Caso Master'
$ws.Range("E80").Value = 'https://www.metropoles.com/colunas/tacio-lorran/em-despacho-ministro-do-tcu-preve-cautelar-contra-bc-no-caso-master'

# Row 81
$ws.Range("A81").Value = '05/01/2026 12:29:54'
$ws.Range("B81").Value = '05/01 12:12'
$ws.Range("C81").Value = 'Metrópoles'
$ws.Range("D81").Value = 'DNA "lixo" pode ter papel no desenvolvimento do Alzheimer, diz estudo'
$ws.Range("E81").Value = 'https://www.metropoles.com/saude/dna-lixo-papel-alzheimer'
$ws.Range("F81").Value = 'lula'
$ws.Range("G81").Value = 'Pesquisa mapeia sinais genéticos em cé&lt;b&gt;lula&lt;/b&gt;s do cérebro e ajuda a entender por que o Alzheimer se desenvolve'

# Row 82
$ws.Range("A82").Value = '05/01/2026 12:29:55'
$ws.Range("B82").Value = '05/01 12:02'
$ws.Range("C82").Value = 'Folha de S.Paulo - Mercado - Principal'
$ws.Range("D82").Value = 'Governo regulamenta corte de 10% nos benefícios fiscais'
$ws.Range("E82").Value = 'https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/blogs/que-imposto-e-esse/2026/01/governo-regulamenta-corte-de-10-nos-beneficios-fiscais.shtml'
$ws.Range("F82").Value = 'senado'
$ws.Range("G82").Value = 'FB nº 2.305&lt;/a&gt;, que regulamenta a &lt;a href="https://www1.folha.uol.com.br/mercado/2025/12/&lt;b&gt;senado&lt;/b&gt;-aprova-corte-de-beneficios-fiscais-e-livra-governo-de-cortar-r-20-bi-no-orcamento.shtml"&gt;'

# Row 83 (no Publicado_em_BRT)
$ws.Range("A83").Value = '05/01/2026 12:29:56'
$ws.Range("C83").Value = 'VEJA'
$ws.Range("D83").Value = 'Para lembrar de invasões golpistas, PT convoca atos em 8 de janeiro'
$ws.Range("E83").Value = 'https://veja.abril.com.br/coluna/radar/para-lembrar-de-invasoes-golpistas-pt-convoca-atos-em-8-de-janeiro/'
$ws.Range("F83").Value = 'lula'
$ws.Range("G83").Value = 'Perfis oficiais ligados ao partido de Lula estão promovendo a mobilização nas ruas nesta quinta-feira'

# Row 84
$ws.Range("A84").Value = '05/01/2026 12:30:00'
$ws.Range("C84").Value = 'VEJA'
$ws.Range("D84").Value = 'Ibovespa abre em alta enquanto crise na Venezuela eleva cautela global'
$ws.Range("E84").Value = 'https://veja.abril.com.br/economia/ibovespa-abre-em-alta-enquanto-crise-na-venezuela-eleva-cautela-global/'
$ws.Range("F84").Value = 'inflação'
$ws.Range("G84").Value = 'Mercado acompanha queda nas projeções de &lt;b&gt;inflação&lt;/b&gt;, tensões geopolíticas envolvendo os EUA e a agenda de indicadores americanos ao longo da '
